# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (a copy of "2021-Q4", so it keeps the
#    same header/column styling) right before the "总计" sheet, and replace
#    its data with the 2022-Q1 fund holdings.
# 2. Insert a new top data-row in "总计" for "2022-Q1" (8 funds, 0.95亿元),
#    shifting the existing history down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by duplicating "2021-Q4" (index 5),
# inserting the copy right before "总计" (index 6), then clearing out the
# rows it doesn't need and overwriting the remaining ones.
# ---------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item(5)
$zongji = $wb.Worksheets.Item(6)
$q4_2021.Copy($zongji)

$q1_2022 = $wb.Worksheets.Item(6)
$q1_2022.Name = "2022-Q1"

# the source sheet had 23 rows of data; only 9 are needed here
$q1_2022.Rows("10:23").Delete()

# header row
$q1_2022.Range("B1").Value = "基金代码"
$q1_2022.Range("C1").Value = "基金名称"
$q1_2022.Range("D1").Value = "基金规模"
$q1_2022.Range("E1").Value = "股票总仓位"
$q1_2022.Range("F1").Value = "仓位占比"
$q1_2022.Range("G1").Value = "持有市值(亿元)"
$q1_2022.Range("H1").Value = "仓位排名"

# data rows - B:G are stored as text, so force text format before writing
# (cleared again afterwards so no stray NumberFormat style sticks to the
# cells, matching the plain/unstyled look of the source data rows)
$q1_2022.Range("B2:G9").NumberFormat = "@"

$q1_2022.Range("A2").Value = 0
$q1_2022.Range("B2").Value = "010755"
$q1_2022.Range("C2").Value = "博道睿见一年持有期混合"
$q1_2022.Range("D2").Value = "6.51"
$q1_2022.Range("E2").Value = "93.26"
$q1_2022.Range("F2").Value = "8.32"
$q1_2022.Range("G2").Value = "0.5416"
$q1_2022.Range("H2").Value = 6

$q1_2022.Range("A3").Value = 1
$q1_2022.Range("B3").Value = "012027"
$q1_2022.Range("C3").Value = "光大保德信安阳一年持有期混合型证券投资基金A"
$q1_2022.Range("D3").Value = "15.22"
$q1_2022.Range("E3").Value = "22.05"
$q1_2022.Range("F3").Value = "1.15"
$q1_2022.Range("G3").Value = "0.1750"
$q1_2022.Range("H3").Value = 2

$q1_2022.Range("A4").Value = 2
$q1_2022.Range("B4").Value = "012028"
$q1_2022.Range("C4").Value = "光大保德信安阳一年持有期混合型证券投资基金C"
$q1_2022.Range("D4").Value = "7.68"
$q1_2022.Range("E4").Value = "22.05"
$q1_2022.Range("F4").Value = "1.15"
$q1_2022.Range("G4").Value = "0.0883"
$q1_2022.Range("H4").Value = 2

$q1_2022.Range("A5").Value = 3
$q1_2022.Range("B5").Value = "011917"
$q1_2022.Range("C5").Value = "山西证券品质生活混合型证券投资基金A"
$q1_2022.Range("D5").Value = "2.16"
$q1_2022.Range("E5").Value = "73.96"
$q1_2022.Range("F5").Value = "3.88"
$q1_2022.Range("G5").Value = "0.0838"
$q1_2022.Range("H5").Value = 6

$q1_2022.Range("A6").Value = 4
$q1_2022.Range("B6").Value = "003848"
$q1_2022.Range("C6").Value = "中银广利灵活配置混合A"
$q1_2022.Range("D6").Value = "5.36"
$q1_2022.Range("E6").Value = "23.61"
$q1_2022.Range("F6").Value = "0.55"
$q1_2022.Range("G6").Value = "0.0295"
$q1_2022.Range("H6").Value = 10

$q1_2022.Range("A7").Value = 5
$q1_2022.Range("B7").Value = "005226"
$q1_2022.Range("C7").Value = "山西证券改革精选灵活配置混合"
$q1_2022.Range("D7").Value = "0.57"
$q1_2022.Range("E7").Value = "88.66"
$q1_2022.Range("F7").Value = "3.99"
$q1_2022.Range("G7").Value = "0.0227"
$q1_2022.Range("H7").Value = 9

$q1_2022.Range("A8").Value = 6
$q1_2022.Range("B8").Value = "011918"
$q1_2022.Range("C8").Value = "山西证券品质生活混合型证券投资基金C"
$q1_2022.Range("D8").Value = "0.21"
$q1_2022.Range("E8").Value = "73.96"
$q1_2022.Range("F8").Value = "3.88"
$q1_2022.Range("G8").Value = "0.0081"
$q1_2022.Range("H8").Value = 6

$q1_2022.Range("A9").Value = 7
$q1_2022.Range("B9").Value = "003849"
$q1_2022.Range("C9").Value = "中银广利灵活配置混合C"
$q1_2022.Range("D9").Value = "0.14"
$q1_2022.Range("E9").Value = "23.61"
$q1_2022.Range("F9").Value = "0.55"
$q1_2022.Range("G9").Value = "0.0008"
$q1_2022.Range("H9").Value = 10

# drop the temporary text NumberFormat now that the values are committed
$q1_2022.Range("B2:G9").ClearFormats()

# ---------------------------------------------------------------------
# Step 2: insert a new top row into "总计" for the 2022-Q1 summary,
# pushing the existing rows down by one.
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$zongji.Rows("2:2").Insert()

# the freshly inserted row has no formatting of its own yet; match the
# plain (unstyled) look of the other data rows for B:D ...
$zongji.Range("B2:D2").ClearFormats()
# ... and the bold/centered/bordered look of column A's index cells
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 8
$zongji.Range("D2").Value = 0.95

# the pre-existing rows kept their old index values (0..4) after the
# shift; renumber column A so the running index stays 0,1,2,3,4,5
$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2
$zongji.Range("A5").Value = 3
$zongji.Range("A6").Value = 4
$zongji.Range("A7").Value = 5

# restore the originally-active sheet/tab (inserting sheets above made the
# newest one active)
$wb.Worksheets.Item(1).Activate()
